$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in column H, matching the style used by the other
# header cells in row 1 (bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-13 (plain numeric, no special style)
$saveValues = @(1,0,0,0,1,0,0,0,0,1,0,0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
